$d = $word.ActiveDocument

# Occurrence 1 (attribution paragraph): "APGLV3 License (see LICENSE file)"
# becomes "AGPL-3.0 license (see LICENSE file)"
$d.Content.Find.Execute("APGLV3 License (see LICENSE file)", $true, $false, $false, $false, $false, $true, 1, $false, "AGPL-3.0 license (see LICENSE file)", 2)

# Occurrence 2 (Open Access Discovery Programme paragraph): "APGLV3" -> "AGPL-3.0 license"
$d.Content.Find.Execute("APGLV3", $true, $false, $false, $false, $false, $true, 1, $false, "AGPL-3.0 license", 2)
